$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text without altering its stored style,
# even when the text looks like a pure number (e.g. "5.360", "1.002").
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = '22.419.79'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.567.91'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +0.00%  '
Set-TextValue "D6" '285.17'
$ws.Range("E6").Value = '  -2.25%  '
Set-TextValue "D7" '0.3619'
$ws.Range("E7").Value = '  -2.80%  '
Set-TextValue "D8" '48.56'
$ws.Range("E8").Value = '  -2.75%  '
Set-TextValue "D9" '0.3324'
$ws.Range("E9").Value = '  -1.91%  '
Set-TextValue "D10" '1.122'
$ws.Range("E10").Value = '  -1.99%  '
Set-TextValue "D11" '0.07389'
$ws.Range("E11").Value = '  -2.34%  '
Set-TextValue "D13" '20.76'
$ws.Range("E13").Value = '  -2.28%  '
Set-TextValue "D14" '5.942'
$ws.Range("E14").Value = '  -1.18%  '
Set-TextValue "D15" '6.902'
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '1.568.99'
$ws.Range("E16").Value = '  -0.40%  '
Set-TextValue "D17" '0.00001103'
$ws.Range("E17").Value = '  -1.64%  '
Set-TextValue "D18" '87.97'
$ws.Range("E18").Value = '  -3.35%  '
Set-TextValue "D19" '0.06696'
$ws.Range("E19").Value = '  -0.87%  '
Set-TextValue "D20" '1.002'
$ws.Range("E20").Value = '  +0.05%  '
Set-TextValue "D21" '6.330'
$ws.Range("E21").Value = '  +0.41%  '
Set-TextValue "D22" '16.18'
$ws.Range("E22").Value = '  -0.68%  '
Set-TextValue "D23" '12.00'
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("D24").Value = '22.420.89'
$ws.Range("E24").Value = '  -0.07%  '
Set-TextValue "D25" '2.376'
$ws.Range("E25").Value = '  +1.71%  '
Set-TextValue "D26" '2.538'
$ws.Range("E26").Value = '  -5.62%  '
Set-TextValue "D27" '150.38'
$ws.Range("E27").Value = '  +1.14%  '
Set-TextValue "D28" '19.40'
$ws.Range("E28").Value = '  -3.44%  '
Set-TextValue "D29" '4.996'
$ws.Range("E29").Value = '  -0.76%  '
Set-TextValue "D30" '123.74'
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("D31").Value = '1.744.75'
$ws.Range("E31").Value = '  -0.34%  '
Set-TextValue "D32" '1.039'
$ws.Range("E32").Value = '  -2.02%  '
Set-TextValue "D33" '2.003'
$ws.Range("E33").Value = '  +0.71%  '
Set-TextValue "D34" '6.084'
$ws.Range("E34").Value = '  -1.21%  '
Set-TextValue "D35" '9.805'
$ws.Range("E35").Value = '  -0.30%  '
Set-TextValue "D36" '0.08236'
$ws.Range("E36").Value = '  -1.44%  '
Set-TextValue "D37" '0.02413'
$ws.Range("E37").Value = '  -2.59%  '
$ws.Range("E38").Value = '  -3.12%  '
Set-TextValue "D39" '0.06409'
$ws.Range("E39").Value = '  -1.71%  '
Set-TextValue "D40" '5.360'
$ws.Range("E40").Value = '  -1.86%  '
Set-TextValue "D41" '1.285'
$ws.Range("E41").Value = '  -5.35%  '
Set-TextValue "D42" '0.6255'
$ws.Range("E42").Value = '  +0.68%  '
Set-TextValue "D43" '11.18'
$ws.Range("E43").Value = '  -1.14%  '
Set-TextValue "D44" '1.002'
$ws.Range("E44").Value = '  +0.03%  '
Set-TextValue "D45" '13.85'
$ws.Range("E45").Value = '  -1.02%  '
Set-TextValue "D46" '0.6049'
$ws.Range("E46").Value = '  +4.17%  '
Set-TextValue "D47" '3.750'
$ws.Range("E47").Value = '  -1.70%  '
Set-TextValue "D48" '2.029'
$ws.Range("E48").Value = '  -1.70%  '
Set-TextValue "D49" '123.32'
$ws.Range("E49").Value = '  -5.18%  '
Set-TextValue "D50" '1.211'
$ws.Range("E50").Value = '  -0.80%  '
Set-TextValue "D51" '0.07207'
$ws.Range("E51").Value = '  -1.53%  '
